# Vignola report: append rolling-7-day new-case data through 2022-01-05
# (commit message: "aggiornamento fino a 6 gennaio 2022")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A = date (serial), B = nuovi pos., C = somma mobile 7gg.,
# D = somma mobile 7gg. per 100mila abitanti
$newRows = @(
    @(465, 44539, 14, 87, 340.522133938706),
    @(466, 44540, 20, 96, 375.7485615875377),
    @(467, 44541, 2, 90, 352.2642764883166),
    @(468, 44542, 18, 90, 352.2642764883166),
    @(469, 44543, 20, 95, 371.8345140710009),
    @(470, 44544, 16, 101, 395.3187991702219),
    @(471, 44545, 5, 95, 371.8345140710009),
    @(472, 44546, 15, 96, 375.7485615875377),
    @(473, 44547, 20, 96, 375.7485615875377),
    @(474, 44548, 4, 98, 383.5766566206113),
    @(475, 44550, 15, 95, 371.8345140710009),
    @(476, 44551, 5, 80, 313.1238013229481),
    @(477, 44552, 4, 68, 266.1552311245059),
    @(478, 44553, 13, 76, 297.4676112568006),
    @(479, 44554, 21, 82, 320.9518963560217),
    @(480, 44555, 12, 74, 289.639516223727),
    @(481, 44556, 20, 90, 352.2642764883166),
    @(482, 44557, 13, 88, 344.4361814552428),
    @(483, 44558, 2, 85, 332.6940389056323),
    @(484, 44559, 13, 94, 367.9204665544639),
    @(485, 44560, 38, 119, 465.7716544678852),
    @(486, 44561, 29, 127, 497.0840346001801),
    @(487, 44562, 46, 161, 630.161650162433),
    @(488, 44563, 38, 179, 700.6145054600963),
    @(489, 44564, 14, 180, 704.5285529766331),
    @(490, 44565, 8, 186, 728.0128380758542),
    @(491, 44566, 30, 203, 794.5516458569807),
)

$lastStyledRow = 464
foreach ($row in $newRows) {
    $r = $row[0]

    # Reuse the date-column formatting (border/bold/centered/date numfmt)
    # from the last existing row, then overwrite with the new values.
    $ws.Cells.Item($lastStyledRow, 1).Copy($ws.Cells.Item($r, 1))
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
